# Rename the "Data" label textboxes to "Header"/"Headers" on the IOAM
# encapsulation example slides, and shrink their auto-fit textbox bounds to
# match the new (shorter) caption text.
#
# PowerPoint's Shape.Width/Height are 32-bit (Single) point values, so a
# plain EMU/12700 conversion can truncate to one EMU below the target once
# it round-trips back through PptxGenJS/PowerPoint's internal EMU storage.
# The point values below were chosen (and verified against this runtime) so
# that they land exactly on the target EMU extents from the authored deck:
#   712054 EMU -> 56.06724409448819 pt
#   307777 EMU -> 24.234411239624023 pt
#   861133 EMU -> 67.80575561523438 pt

$p = $ppt.ActivePresentation

$WIDTH_HEADER_PT  = 56.06724409448819    # 712054 EMU
$WIDTH_HEADERS_PT = 67.80575561523438    # 861133 EMU
$HEIGHT_PT        = 24.234411239624023   # 307777 EMU

function Resize-DataLabel($Slide, [string]$ShapeName, [string]$NewText, [double]$NewWidthPt, [double]$NewHeightPt) {
    $shape = $Slide.Shapes.Item($ShapeName)
    $shape.TextFrame.TextRange.Text = $NewText
    $shape.Width = $NewWidthPt
    $shape.Height = $NewHeightPt
}

# Slide 21: "Generic PW Control Word [RFC4385] with IOAM Data Fields"
Resize-DataLabel $p.Slides.Item(21) "TextBox 9" "Header" $WIDTH_HEADER_PT $HEIGHT_PT

# Slide 22: "MPLS Encap with Additional G-ACh [RFC5586] with IOAM Data Fields"
Resize-DataLabel $p.Slides.Item(22) "TextBox 9" "Header" $WIDTH_HEADER_PT $HEIGHT_PT

# Slide 26: "Example - Generic Delivery Function with IOAM Data Fields"
Resize-DataLabel $p.Slides.Item(26) "TextBox 10" "Headers" $WIDTH_HEADERS_PT $HEIGHT_PT

# Slide 27: "Example - Generic Delivery Function with IOAM Data Fields and PW"
Resize-DataLabel $p.Slides.Item(27) "TextBox 9" "Headers" $WIDTH_HEADERS_PT $HEIGHT_PT
